# SSCS hearing list update
# - Rename the "Panel" column to "Tribunal"
# - Replace the single "Panel 1" value with a multi-line list of tribunal
#   members, wrap the text, widen the column and grow the row to fit
# - Move the active selection to I17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Panel" -> "Tribunal"
$ws.Range("G1").Value = "Tribunal"

# Data cell: "Panel 1" -> multi-line tribunal members, with wrap text
$ws.Range("G2").Value = "Tribunal member 1`nTribunal member 2`nTribunal member 3"
$ws.Range("G2").WrapText = $true

# Column G grows to fit the longer wrapped text (stored width "41")
$ws.Columns(7).ColumnWidth = 40.166666666666664

# Row 2 grows to fit the three wrapped lines
$ws.Rows(2).RowHeight = 48

# Move the selection
[void]$ws.Range("I17").Select()
